$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly imprecise timestamp value in A11 (rounding fix)
$ws.Range("A11").Value = 45864.62531725694

# Append new row 12 with the latest automated sensor reading
$ws.Range("A12").Value = 45864.66703143053
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat

$ws.Range("B12").Value = 2025
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 22.22
$ws.Range("E12").Value = 64.34999999999999
$ws.Range("F12").Value = 286.52
$ws.Range("G12").Value = 9.35
$ws.Range("H12").Value = "ESE"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "16:00:31"
